# Apply "6-9 introducing Arduino CAN Shield Functions to STM Libraries" edit
# - Typography sheet: bump Typography_03 size, add Typography_15 / Typography_16
# - Translation sheet: tweak alignment/typography/text on existing rows,
#   drop the old "START MOTOR" row, append a new "Start Engine" row

$wb = $excel.ActiveWorkbook

$typography = $wb.Worksheets.Item("Typography")
$translation = $wb.Worksheets.Item("Translation")

# --- Typography sheet -------------------------------------------------

# Typography_03 (row 10) size 19 -> 22
$typography.Range("D10").Value = 22

# New rows 22 and 23 in the Table7 listing
$typography.Range("B22").Value = "Typography_15"
$typography.Range("C22").Value = "malgunbd.ttf"
$typography.Range("D22").Value = 22
$typography.Range("E22").Value = 4
$typography.Range("F22").Value = "?"

$typography.Range("B23").Value = "Typography_16"
$typography.Range("C23").Value = "malgunbd.ttf"
$typography.Range("D23").Value = 21
$typography.Range("E23").Value = 4
$typography.Range("F23").Value = "?"

# --- Translation sheet -------------------------------------------------

# Row 4 (SingleUseId1 / Default): alignment Center -> Right
$translation.Range("D4").Value = "Right"

# Row 5 (SingleUseId2): now uses new Typography_15 + updated GB text
$translation.Range("C5").Value = "Typography_15"
$translation.Range("F5").Value = "HOLD FOOT ON BRAKE!"

# The old row 19 (SingleUseId25 / START MOTOR) is removed entirely.
$translation.Rows("19:19").Delete()

# The (now shifted-up) row 17 (SingleUseId22 / "<>") is pulled out of the
# middle of the table and relocated to the end of the list.
$translation.Rows("17:17").Delete()

# New row 32, appended after the shifted "Startup Screen" row (31)
$translation.Range("B32").Value = "SingleUseId40"
$translation.Range("C32").Value = "Default"
$translation.Range("D32").Value = "Center"
$translation.Range("E32").Value = "LTR"
$translation.Range("F32").Value = "Start Engine"

# Row 33: the relocated SingleUseId22 row
$translation.Range("B33").Value = "SingleUseId22"
$translation.Range("C33").Value = "Typography_10"
$translation.Range("D33").Value = "Center"
$translation.Range("E33").Value = "LTR"
$translation.Range("F33").Value = "<>"
